$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.445154666666667
$ws.Range("H2").Value = 10.335464
$ws.Range("I2").Value = 0.01110365039942287
$ws.Range("J2").Value = 0.01110365039942286
$ws.Range("M2").Value = 10.34761366666667
$ws.Range("N2").Value = 31.042841
$ws.Range("O2").Value = 0.2299953477621856
$ws.Range("P2").Value = 0.2299953477621856
$ws.Range("Q2").Value = 35.64912951258044
$ws.Range("R2").Value = 320.842165613224
$ws.Range("S2").Value = 0.002553787935044992
$ws.Range("T2").Value = 0.002553787935044993
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.445154666666667
$ws.Range("H3").Value = 10.335464
$ws.Range("I3").Value = 0.01110365039942287
$ws.Range("J3").Value = 0.01110365039942286
$ws.Range("O3").Value = 0.6794731949692173
$ws.Range("P3").Value = 0.6794731949692174
$ws.Range("Q3").Value = 105.3179038770409
$ws.Range("R3").Value = 947.861134893368
$ws.Range("S3").Value = 0.007544632812717081
$ws.Range("T3").Value = 0.007544632812717081
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.445154666666667
$ws.Range("H4").Value = 10.335464
$ws.Range("I4").Value = 0.01110365039942287
$ws.Range("J4").Value = 0.01110365039942286
$ws.Range("M4").Value = 4.073058666666666
$ws.Range("N4").Value = 12.219176
$ws.Range("O4").Value = 0.09053145726859702
$ws.Range("P4").Value = 0.09053145726859703
$ws.Range("Q4").Value = 14.03231707307378
$ws.Range("R4").Value = 126.290853657664
$ws.Range("S4").Value = 0.001005229651660791
$ws.Range("T4").Value = 0.001005229651660791
# Row 5
$ws.Range("I5").Value = 0.938949437922138
$ws.Range("J5").Value = 0.938949437922138
$ws.Range("M5").Value = 10.34761366666667
$ws.Range("N5").Value = 31.042841
$ws.Range("O5").Value = 0.2299953477621856
$ws.Range("P5").Value = 0.2299953477621856
$ws.Range("Q5").Value = 3014.569885953066
$ws.Range("R5").Value = 27131.1289735776
$ws.Range("S5").Value = 0.2159540025060108
$ws.Range("T5").Value = 0.2159540025060109
# Row 6
$ws.Range("I6").Value = 0.938949437922138
$ws.Range("J6").Value = 0.938949437922138
$ws.Range("O6").Value = 0.6794731949692173
$ws.Range("P6").Value = 0.6794731949692174
$ws.Range("S6").Value = 0.6379909744995059
$ws.Range("T6").Value = 0.637990974499506
# Row 7
$ws.Range("I7").Value = 0.938949437922138
$ws.Range("J7").Value = 0.938949437922138
$ws.Range("M7").Value = 4.073058666666666
$ws.Range("N7").Value = 12.219176
$ws.Range("O7").Value = 0.09053145726859702
$ws.Range("P7").Value = 0.09053145726859703
$ws.Range("Q7").Value = 1186.604022510712
$ws.Range("R7").Value = 10679.43620259641
$ws.Range("S7").Value = 0.08500446091662123
$ws.Range("T7").Value = 0.08500446091662124
# Row 8
$ws.Range("G8").Value = 15.497141
$ws.Range("H8").Value = 46.491423
$ws.Range("I8").Value = 0.04994691167843914
$ws.Range("J8").Value = 0.04994691167843914
$ws.Range("M8").Value = 10.34761366666667
$ws.Range("N8").Value = 31.042841
$ws.Range("O8").Value = 0.2299953477621856
$ws.Range("P8").Value = 0.2299953477621856
$ws.Range("Q8").Value = 160.3584280058603
$ws.Range("R8").Value = 1443.225852052743
$ws.Range("S8").Value = 0.01148755732112978
$ws.Range("T8").Value = 0.01148755732112978
# Row 9
$ws.Range("G9").Value = 15.497141
$ws.Range("H9").Value = 46.491423
$ws.Range("I9").Value = 0.04994691167843914
$ws.Range("J9").Value = 0.04994691167843914
$ws.Range("O9").Value = 0.6794731949692173
$ws.Range("P9").Value = 0.6794731949692174
$ws.Range("Q9").Value = 473.7454669302556
$ws.Range("R9").Value = 4263.709202372301
$ws.Range("S9").Value = 0.03393758765699436
$ws.Range("T9").Value = 0.03393758765699436
# Row 10
$ws.Range("G10").Value = 15.497141
$ws.Range("H10").Value = 46.491423
$ws.Range("I10").Value = 0.04994691167843914
$ws.Range("J10").Value = 0.04994691167843914
$ws.Range("M10").Value = 4.073058666666666
$ws.Range("N10").Value = 12.219176
$ws.Range("O10").Value = 0.09053145726859702
$ws.Range("P10").Value = 0.09053145726859703
$ws.Range("Q10").Value = 63.12076445860532
$ws.Range("R10").Value = 568.0868801274479
$ws.Range("S10").Value = 0.004521766700315002
$ws.Range("T10").Value = 0.004521766700315003
